$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053970280515662
$ws.Range("D2").Value = 1.060183157554081
$ws.Range("E2").Value = 1.060729293787088
$ws.Range("F2").Value = 1.071783102435846
$ws.Range("I2").Value = 1.047391627475972
$ws.Range("J2").Value = 1.058984830962287
$ws.Range("K2").Value = 1.062910631262459
$ws.Range("L2").Value = 1.06345528113359
$ws.Range("M2").Value = 1.074479357892147
$ws.Range("N2").Value = 1.023404877927164
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.054982471454197
$ws.Range("D3").Value = 1.061000159801093
$ws.Range("E3").Value = 1.061628642693826
$ws.Range("F3").Value = 1.072760380921049
$ws.Range("I3").Value = 1.047665613840825
$ws.Range("J3").Value = 1.059647812783776
$ws.Range("K3").Value = 1.063542002362599
$ws.Range("L3").Value = 1.064168898251827
$ws.Range("M3").Value = 1.075272859041089
$ws.Range("N3").Value = 1.023630156382624
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.055637854959629
$ws.Range("D4").Value = 1.061529188587655
$ws.Range("E4").Value = 1.062211311221711
$ws.Range("F4").Value = 1.073393598220148
$ws.Range("I4").Value = 1.047841928800197
$ws.Range("J4").Value = 1.060076612137426
$ws.Range("K4").Value = 1.063950261889957
$ws.Range("L4").Value = 1.064630748668571
$ws.Range("M4").Value = 1.075786536574731
$ws.Range("N4").Value = 1.023775752224478
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.055913480111331
$ws.Range("D5").Value = 1.061751680689362
$ws.Range("E5").Value = 1.062456438694323
$ws.Range("F5").Value = 1.073660005845333
$ws.Range("I5").Value = 1.0479158183614
$ws.Range("J5").Value = 1.060256832158485
$ws.Range("K5").Value = 1.064121826429226
$ws.Range("L5").Value = 1.064824931879077
$ws.Range("M5").Value = 1.076002540567603
$ws.Range("N5").Value = 1.023836918627022
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.055959764746618
$ws.Range("D6").Value = 1.061789043245141
$ws.Range("E6").Value = 1.062497606826318
$ws.Range("F6").Value = 1.073704748734734
$ws.Range("I6").Value = 1.047928211053954
$ws.Range("J6").Value = 1.06028708910241
$ws.Range("K6").Value = 1.064150628873753
$ws.Range("L6").Value = 1.064857537319235
$ws.Range("M6").Value = 1.076038811725881
$ws.Range("N6").Value = 1.023847186257088
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.055641537476639
$ws.Range("D7").Value = 1.061532161191987
$ws.Range("E7").Value = 1.062214585946094
$ws.Range("F7").Value = 1.073397157175287
$ws.Range("I7").Value = 1.047842917033152
$ws.Range("J7").Value = 1.060079020432648
$ws.Range("K7").Value = 1.063952554610976
$ws.Range("L7").Value = 1.064633343271493
$ws.Range("M7").Value = 1.075789422619786
$ws.Range("N7").Value = 1.023776569698752
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.054312266137109
$ws.Range("D8").Value = 1.060459189270613
$ws.Range("E8").Value = 1.061033081247354
$ws.Range("F8").Value = 1.072113201052326
$ws.Range("I8").Value = 1.04748442356844
$ws.Range("J8").Value = 1.059208928546532
$ws.Range("K8").Value = 1.063124063543415
$ws.Range("L8").Value = 1.063696431742859
$ws.Range("M8").Value = 1.074747477302838
$ws.Range("N8").Value = 1.023481047737489
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.051973227055887
$ws.Range("D9").Value = 1.058571384730558
$ws.Range("E9").Value = 1.058956753798304
$ws.Range("F9").Value = 1.069857285270036
$ws.Range("I9").Value = 1.046845282341331
$ws.Range("J9").Value = 1.05767426024382
$ws.Range("K9").Value = 1.061662048084429
$ws.Range("L9").Value = 1.062046218417307
$ws.Range("M9").Value = 1.072913232761911
$ws.Range("N9").Value = 1.022958979780124
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.050416125161905
$ws.Range("D10").Value = 1.05731486654599
$ws.Range("E10").Value = 1.057576385648147
$ws.Range("F10").Value = 1.068357827866676
$ws.Range("I10").Value = 1.046414217169599
$ws.Range("J10").Value = 1.056650211998467
$ws.Range("K10").Value = 1.060685999146492
$ws.Range("L10").Value = 1.06094662510515
$ws.Range("M10").Value = 1.071691664371541
$ws.Range("N10").Value = 1.022610066591584
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049742424886378
$ws.Range("D11").Value = 1.05677127177604
$ws.Range("E11").Value = 1.056979596834605
$ws.Range("F11").Value = 1.067709623294419
$ws.Range("I11").Value = 1.04622638708637
$ws.Range("J11").Value = 1.056206574354415
$ws.Range("K11").Value = 1.060263043044137
$ws.Range("L11").Value = 1.060470629101546
$ws.Range("M11").Value = 1.071163021873343
$ws.Range("N11").Value = 1.022458781566475
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049492263456384
$ws.Range("D12").Value = 1.056569430464538
$ws.Range("E12").Value = 1.056758062048863
$ws.Range("F12").Value = 1.067469013062322
$ws.Range("I12").Value = 1.046156442307366
$ws.Range("J12").Value = 1.05604175546833
$ws.Range("K12").Value = 1.060105890867979
$ws.Range("L12").Value = 1.06029384404568
$ws.Range("M12").Value = 1.07096670711974
$ws.Range("N12").Value = 1.022402557338941
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04954592026639
$ws.Range("D13").Value = 1.056612722746936
$ws.Range("E13").Value = 1.05680557569982
$ws.Range("F13").Value = 1.067520617435749
$ws.Range("I13").Value = 1.04617145367306
$ws.Range("J13").Value = 1.05607711112173
$ws.Range("K13").Value = 1.060139602679075
$ws.Range("L13").Value = 1.060331764075522
$ws.Range("M13").Value = 1.071008815178265
$ws.Range("N13").Value = 1.022414618989709
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049721744812598
$ws.Range("D14").Value = 1.056754585995173
$ws.Range("E14").Value = 1.056961281855801
$ws.Range("F14").Value = 1.067689731070988
$ws.Range("I14").Value = 1.046220609023955
$ws.Range("J14").Value = 1.056192951025243
$ws.Range("K14").Value = 1.060250053760093
$ws.Range("L14").Value = 1.060456015559858
$ws.Range("M14").Value = 1.071146793478945
$ws.Range("N14").Value = 1.02245413466923
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.049830086845388
$ws.Range("D15").Value = 1.056842002458454
$ws.Range("E15").Value = 1.05705723602564
$ws.Range("F15").Value = 1.067793949033976
$ws.Range("I15").Value = 1.046250871904314
$ws.Range("J15").Value = 1.056264319577012
$ws.Range("K15").Value = 1.060318100050159
$ws.Range("L15").Value = 1.060532573812395
$ws.Range("M15").Value = 1.071231812664867
$ws.Range("N15").Value = 1.02247847758671
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.050460847743304
$ws.Range("D16").Value = 1.057350953462256
$ws.Range("E16").Value = 1.057616012016641
$ws.Range("F16").Value = 1.068400869705843
$ws.Range("I16").Value = 1.046426658086957
$ws.Range("J16").Value = 1.05667965022939
$ws.Range("K16").Value = 1.060714062657721
$ws.Range("L16").Value = 1.060978218303626
$ws.Range("M16").Value = 1.071726755123394
$ws.Range("N16").Value = 1.022620102623535
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.05085665097487
$ws.Range("D17").Value = 1.057670335766647
$ws.Range("E17").Value = 1.057966764352869
$ws.Range("F17").Value = 1.068781862281026
$ws.Range("I17").Value = 1.046536609515034
$ws.Range("J17").Value = 1.056940118261418
$ws.Range("K17").Value = 1.060962354234045
$ws.Range("L17").Value = 1.061257796031699
$ws.Range("M17").Value = 1.072037301609721
$ws.Range("N17").Value = 1.022708886055279
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.051087567903972
$ws.Range("D18").Value = 1.05785667287675
$ws.Range("E18").Value = 1.058171440928901
$ws.Range("F18").Value = 1.069004192013359
$ws.Range("I18").Value = 1.046600628829642
$ws.Range("J18").Value = 1.057092023781227
$ws.Range("K18").Value = 1.061107147541863
$ws.Range("L18").Value = 1.061420881980791
$ws.Range("M18").Value = 1.072218467461908
$ws.Range("N18").Value = 1.022760652291314
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051166313348027
$ws.Range("D19").Value = 1.057920216926553
$ws.Range("E19").Value = 1.058241245362738
$ws.Range("F19").Value = 1.069080018182389
$ws.Range("I19").Value = 1.046622438500181
$ws.Range("J19").Value = 1.057143816043544
$ws.Range("K19").Value = 1.061156513044101
$ws.Range("L19").Value = 1.061476492255592
$ws.Range("M19").Value = 1.072280245287127
$ws.Range("N19").Value = 1.022778299899674
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05081417969334
$ws.Range("D20").Value = 1.057636064218467
$ws.Range("E20").Value = 1.057929122767769
$ws.Range("F20").Value = 1.068740974700756
$ws.Range("I20").Value = 1.04652482450388
$ws.Range("J20").Value = 1.056912174682012
$ws.Range("K20").Value = 1.060935718095781
$ws.Range("L20").Value = 1.06122779863828
$ws.Range("M20").Value = 1.07200397986981
$ws.Range("N20").Value = 1.022699362462234
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049669966630558
$ws.Range("D21").Value = 1.056712808753988
$ws.Range("E21").Value = 1.056915426424622
$ws.Range("F21").Value = 1.067639926845778
$ws.Range("I21").Value = 1.046206138862747
$ws.Range("J21").Value = 1.056158839942195
$ws.Range("K21").Value = 1.060217529973809
$ws.Range("L21").Value = 1.060419426007408
$ws.Range("M21").Value = 1.071106161015093
$ws.Range("N21").Value = 1.022442499118591
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048951022459842
$ws.Range("D22").Value = 1.056132749971592
$ws.Range("E22").Value = 1.056278880873505
$ws.Range("F22").Value = 1.066948591020991
$ws.Range("I22").Value = 1.046004748582103
$ws.Range("J22").Value = 1.055685003019493
$ws.Range("K22").Value = 1.059765702427745
$ws.Range("L22").Value = 1.059911292231849
$ws.Range("M22").Value = 1.070541936741418
$ws.Range("N22").Value = 1.022280824174124
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049332103854059
$ws.Range("D23").Value = 1.056440209152786
$ws.Range("E23").Value = 1.056616249028372
$ws.Range("F23").Value = 1.067314992104459
$ws.Range("I23").Value = 1.046111605940395
$ws.Range("J23").Value = 1.055936210303818
$ws.Range("K23").Value = 1.060005250512515
$ws.Range("L23").Value = 1.060180651728701
$ws.Range("M23").Value = 1.070841016825232
$ws.Range("N23").Value = 1.022366547582616
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.050833370477834
$ws.Range("D24").Value = 1.057651549911308
$ws.Range("E24").Value = 1.057946131105896
$ws.Range("F24").Value = 1.068759449719459
$ws.Range("I24").Value = 1.046530149993826
$ws.Range("J24").Value = 1.056924801247872
$ws.Range("K24").Value = 1.060947753915843
$ws.Range("L24").Value = 1.061241353129043
$ws.Range("M24").Value = 1.072019036439146
$ws.Range("N24").Value = 1.022703665824691
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.052577529347698
$ws.Range("D25").Value = 1.059059076302033
$ws.Range("E25").Value = 1.059492860117440
$ws.Range("F25").Value = 1.070439707437744
$ws.Range("I25").Value = 1.047011393645549
$ws.Range("J25").Value = 1.058071176882173
$ws.Range("K25").Value = 1.062040259288155
$ws.Range("L25").Value = 1.062472744917157
$ws.Range("M25").Value = 1.07338721056367
$ws.Range("N25").Value = 1.02309410129297
